$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 115
$ws.Range("F3").Value = 131
$ws.Range("F4").Value = 877
$ws.Range("F5").Value = 1055
$ws.Range("F8").Value = 651
$ws.Range("F9").Value = 11874
$ws.Range("F11").Value = 2133
$ws.Range("F12").Value = 903
$ws.Range("F13").Value = 245
$ws.Range("F16").Value = 1213
$ws.Range("F17").Value = 184
$ws.Range("F18").Value = 259
$ws.Range("F19").Value = 750
$ws.Range("F20").Value = 665
$ws.Range("F21").Value = 283
$ws.Range("F22").Value = 2911
$ws.Range("F23").Value = 744
$ws.Range("F24").Value = 3769
$ws.Range("F25").Value = 3769
$ws.Range("F27").Value = 830
$ws.Range("F31").Value = 1012
$ws.Range("F34").Value = 265
$ws.Range("F37").Value = 20
$ws.Range("F38").Value = 4281
$ws.Range("F40").Value = 4483
$ws.Range("F41").Value = 5509
$ws.Range("F43").Value = 121
$ws.Range("F44").Value = 51
$ws.Range("F45").Value = 167
$ws.Range("F46").Value = 283
$ws.Range("F47").Value = 73
$ws.Range("F49").Value = 4099

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 4166
$ws.Range("F5").Value = 96
$ws.Range("F12").Value = 829

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 759
$ws.Range("F3").Value = 427

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 759
$ws.Range("F3").Value = 427
$ws.Range("F5").Value = 131
$ws.Range("F6").Value = 877
$ws.Range("F7").Value = 1055
$ws.Range("F10").Value = 651
$ws.Range("F11").Value = 11874
$ws.Range("F12").Value = 2133
$ws.Range("F13").Value = 903
$ws.Range("F14").Value = 245
$ws.Range("F15").Value = 1213
$ws.Range("F16").Value = 184
$ws.Range("F17").Value = 259
$ws.Range("F18").Value = 4166
$ws.Range("F19").Value = 750
$ws.Range("F20").Value = 283
$ws.Range("F21").Value = 744
$ws.Range("F22").Value = 3769
$ws.Range("F24").Value = 96
$ws.Range("F26").Value = 830
$ws.Range("F29").Value = 1012
$ws.Range("F32").Value = 265
$ws.Range("F34").Value = 20
$ws.Range("F35").Value = 4483
$ws.Range("F37").Value = 121
$ws.Range("F38").Value = 167
$ws.Range("F39").Value = 283
$ws.Range("F43").Value = 73
$ws.Range("F45").Value = 4099
